# Apply updated cryptocurrency price/volume data to Sheet1.
# Note: a handful of "Price" values are digit strings that look like plain
# numbers (e.g. "1.00", "198.16"). They are written with a leading
# apostrophe so Excel stores them as literal text (matching the source
# data, which keeps these as text cells) instead of silently coercing them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.014.78'
$ws.Range('E2').Value = '  -1.72%  '

$ws.Range('D3').Value = '3.472.71'
$ws.Range('E3').Value = '  -4.43%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = '''198.16'
$ws.Range('E5').Value = '  +1.02%  '

$ws.Range('D6').Value = '''546.75'
$ws.Range('E6').Value = '  -5.03%  '

$ws.Range('D7').Value = '3.466.72'
$ws.Range('E7').Value = '  -4.45%  '

$ws.Range('D8').Value = '''0.602'
$ws.Range('E8').Value = '  -3.04%  '

$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').Value = '''0.646'
$ws.Range('E10').Value = '  -4.70%  '

$ws.Range('D11').Value = '''61.58'
$ws.Range('E11').Value = '  +9.99%  '

$ws.Range('D12').Value = '''0.141'
$ws.Range('E12').Value = '  -7.63%  '

$ws.Range('D13').Value = '''0.0000266'
$ws.Range('E13').Value = '  -9.69%  '

$ws.Range('D14').Value = '''9.67'
$ws.Range('E14').Value = '  -4.36%  '

$ws.Range('D15').Value = '4.046.87'
$ws.Range('E15').Value = '  -3.84%  '

$ws.Range('D16').Value = '3.489.34'
$ws.Range('E16').Value = '  -3.87%  '

$ws.Range('E17').Value = '  -2.20%  '

$ws.Range('D18').Value = '66.765.14'
$ws.Range('E18').Value = '  -2.01%  '

$ws.Range('D19').Value = '''18.08'
$ws.Range('E19').Value = '  -2.44%  '

$ws.Range('D20').Value = '''11.64'
$ws.Range('E20').Value = '  -7.14%  '

$ws.Range('D21').Value = '''1.01'
$ws.Range('E21').Value = '  -6.37%  '

$ws.Range('D22').Value = '''385.08'
$ws.Range('E22').Value = '  -4.44%  '

$ws.Range('D23').Value = '''3.96'
$ws.Range('E23').Value = '  -6.41%  '

$ws.Range('D24').Value = '''11.77'
$ws.Range('E24').Value = '  -7.86%  '

$ws.Range('D25').Value = '''81.73'
$ws.Range('E25').Value = '  -5.09%  '

$ws.Range('D26').Value = '''3.81'
$ws.Range('E26').Value = '  -1.18%  '

$ws.Range('D27').Value = '''12.04'
$ws.Range('E27').Value = '  -4.69%  '

$ws.Range('D28').Value = '''2.76'
$ws.Range('E28').Value = '  -6.36%  '

$ws.Range('D29').Value = '''8.69'
$ws.Range('E29').Value = '  -5.04%  '

$ws.Range('D30').Value = '''30.69'
$ws.Range('E30').Value = '  -3.29%  '

$ws.Range('D31').Value = '''670.56'
$ws.Range('E31').Value = '  -3.74%  '

$ws.Range('D32').Value = '''6.86'
$ws.Range('E32').Value = '  -15.75%  '

$ws.Range('D33').Value = '''11.56'
$ws.Range('E33').Value = '  -5.28%  '

$ws.Range('D34').Value = '''63.35'
$ws.Range('E34').Value = '  -2.19%  '

$ws.Range('E35').Value = '  -8.09%  '

$ws.Range('D36').Value = '''38.07'
$ws.Range('E36').Value = '  -10.93%  '

$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('D38').Value = '''0.392'
$ws.Range('E38').Value = '  -5.76%  '

$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('D40').Value = '3.040.49'
$ws.Range('E40').Value = '  -4.06%  '

$ws.Range('D41').Value = '''0.128'
$ws.Range('E41').Value = '  -5.72%  '

$ws.Range('D42').Value = '''2.94'
$ws.Range('E42').Value = '  -5.91%  '

$ws.Range('D43').Value = '0.0₃0663'
$ws.Range('E43').Value = '  -16.77%  '

$ws.Range('D44').Value = '''2.77'
$ws.Range('E44').Value = '  +6.40%  '

$ws.Range('D45').Value = '''2.46'
$ws.Range('E45').Value = '  -13.99%  '

$ws.Range('D46').Value = '''2.69'
$ws.Range('E46').Value = '  -7.47%  '

$ws.Range('D47').Value = '''0.0391'
$ws.Range('E47').Value = '  -7.67%  '

$ws.Range('D48').Value = '''0.125'
$ws.Range('E48').Value = '  -5.48%  '

$ws.Range('D49').Value = '''136.03'
$ws.Range('E49').Value = '  -4.57%  '

$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '''8.10'
$ws.Range('E50').Value = '  -8.52%  '

$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').Value = '''2.84'
$ws.Range('E51').Value = '  -8.67%  '
